# Readme.docx edit:
#   - Merge the two paragraphs into a single paragraph.
#   - Replace the second paragraph's content with the new explanatory text.
#   - Justify the resulting paragraph (Alignment = wdAlignParagraphJustify).

$d = $word.ActiveDocument

$p1 = $d.Paragraphs.Item(1)

# Insert the new tail text right before paragraph 1's own paragraph mark, so it
# becomes part of paragraph 1's content (kept as a plain insertion, not a
# replace-over-existing-text edit, so it lands in its own run).
$insertPoint = $d.Range($p1.Range.End - 1, $p1.Range.End - 1)
$insertPoint.InsertBefore(" Within each folder the files include the necessary documents (Parameter values, and initial conditions) for reproducing corresponding figure in the main text")

# The old second paragraph (with the superseded "Each folder Figure i, ..."
# text) now follows; remove it completely, mark and all, which merges the
# remaining content into a single paragraph.
$p2 = $d.Paragraphs.Item(2)
$p2.Range.Delete()

# Justify the merged paragraph.
$d.Paragraphs.Item(1).Alignment = 3
